$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: new entry (UDP, 66.22.244.139)
$ws.Range("A2").Value = "66.22.244.139"
$ws.Range("B2").Value = 50004
$ws.Range("C2").Value = "Ether / IP / UDP 192.168.1.80:51280 > 66.22.244.139:50004 / Raw"

# Row 3: previously row 2's data (162.159.135.234)
$ws.Range("A3").Value = "162.159.135.234"
$ws.Range("B3").Value = 443
$ws.Range("C3").Value = "Ether / IP / TCP 192.168.1.80:52490 > 162.159.135.234:https A"

# Row 4: new entry (146.75.117.44)
$ws.Range("A4").Value = "146.75.117.44"
$ws.Range("B4").Value = 443
$ws.Range("C4").Value = "Ether / IP / TCP 192.168.1.80:53984 > 146.75.117.44:https A / Raw"

# Remove row 5 entirely (delete the whole row, shifting cells up)
$ws.Rows("5").Delete()
